$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 9
$ws.Range("D2").Value = 44509
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19500
$ws.Range("Q2").Value = "$/bandeja 8 kilos"
$ws.Range("S2").Value = 2438
$ws.Range("T2").Value = 8

# Row 3 <- original row 15
$ws.Range("D3").Value = 44159
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 300
$ws.Range("N3").Value = 2000
$ws.Range("O3").Value = 2100
$ws.Range("P3").Value = 2050
$ws.Range("Q3").Value = "$/kilo (en caja de 14 kilos)"
$ws.Range("S3").Value = 2050
$ws.Range("T3").Value = 1

# Row 4 <- original row 2
$ws.Range("D4").Value = 44512
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 19000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 19500
$ws.Range("Q4").Value = "$/bandeja 8 kilos"
$ws.Range("S4").Value = 2438
$ws.Range("T4").Value = 8

# Row 5 <- original row 3
$ws.Range("D5").Value = 44495
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 270
$ws.Range("N5").Value = 19000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 19556
$ws.Range("Q5").Value = "$/bandeja 8 kilos"
$ws.Range("S5").Value = 2444
$ws.Range("T5").Value = 8

# Row 6 <- original row 11
$ws.Range("D6").Value = 44488
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 160
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 17500
$ws.Range("Q6").Value = "$/bandeja 8 kilos"
$ws.Range("S6").Value = 2188
$ws.Range("T6").Value = 8

# Row 7 <- original row 19
$ws.Range("D7").Value = 44519
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 400
$ws.Range("N7").Value = 21000
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 21500
$ws.Range("Q7").Value = "$/bandeja 8 kilos"
$ws.Range("S7").Value = 2688
$ws.Range("T7").Value = 8

# Row 8 <- original row 20
$ws.Range("D8").Value = 44519
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 18000
$ws.Range("Q8").Value = "$/bandeja 8 kilos"
$ws.Range("S8").Value = 2250
$ws.Range("T8").Value = 8

# Row 9 <- original row 16
$ws.Range("D9").Value = 44533
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 18000
$ws.Range("O9").Value = 19000
$ws.Range("P9").Value = 18500
$ws.Range("Q9").Value = "$/bandeja 8 kilos"
$ws.Range("S9").Value = 2312
$ws.Range("T9").Value = 8

# Row 10 <- original row 17
$ws.Range("D10").Value = 44533
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 16000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 16000
$ws.Range("Q10").Value = "$/bandeja 8 kilos"
$ws.Range("S10").Value = 2000
$ws.Range("T10").Value = 8

# Row 11 <- original row 13
$ws.Range("D11").Value = 44491
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 18000
$ws.Range("O11").Value = 19000
$ws.Range("P11").Value = 18500
$ws.Range("Q11").Value = "$/bandeja 8 kilos"
$ws.Range("S11").Value = 2312
$ws.Range("T11").Value = 8

# Row 12 <- original row 4
$ws.Range("D12").Value = 44523
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 400
$ws.Range("N12").Value = 21000
$ws.Range("O12").Value = 22000
$ws.Range("P12").Value = 21500
$ws.Range("Q12").Value = "$/bandeja 8 kilos"
$ws.Range("S12").Value = 2688
$ws.Range("T12").Value = 8

# Row 13 <- original row 5
$ws.Range("D13").Value = 44523
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 18000
$ws.Range("O13").Value = 18000
$ws.Range("P13").Value = 18000
$ws.Range("Q13").Value = "$/bandeja 8 kilos"
$ws.Range("S13").Value = 2250
$ws.Range("T13").Value = 8

# Row 14 <- original row 10
$ws.Range("D14").Value = 44526
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 21000
$ws.Range("O14").Value = 21000
$ws.Range("P14").Value = 21000
$ws.Range("Q14").Value = "$/bandeja 8 kilos"
$ws.Range("S14").Value = 2625
$ws.Range("T14").Value = 8

# Row 15 <- original row 12
$ws.Range("D15").Value = 44516
$ws.Range("L15").Value = "Segunda"
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 18000
$ws.Range("O15").Value = 19000
$ws.Range("P15").Value = 18500
$ws.Range("Q15").Value = "$/bandeja 8 kilos"
$ws.Range("S15").Value = 2312
$ws.Range("T15").Value = 8

# Row 16 <- original row 14
$ws.Range("D16").Value = 44162
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 2000
$ws.Range("O16").Value = 2100
$ws.Range("P16").Value = 2050
$ws.Range("Q16").Value = "$/kilo (en caja de 14 kilos)"
$ws.Range("S16").Value = 2050
$ws.Range("T16").Value = 1

# Row 17 <- original row 7
$ws.Range("D17").Value = 44530
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 19000
$ws.Range("O17").Value = 20000
$ws.Range("P17").Value = 19500
$ws.Range("Q17").Value = "$/bandeja 8 kilos"
$ws.Range("S17").Value = 2438
$ws.Range("T17").Value = 8

# Row 18 <- original row 8
$ws.Range("D18").Value = 44530
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 16000
$ws.Range("O18").Value = 16000
$ws.Range("P18").Value = 16000
$ws.Range("Q18").Value = "$/bandeja 8 kilos"
$ws.Range("S18").Value = 2000
$ws.Range("T18").Value = 8

# Row 19 <- original row 6
$ws.Range("D19").Value = 44498
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 300
$ws.Range("N19").Value = 19000
$ws.Range("O19").Value = 20000
$ws.Range("P19").Value = 19500
$ws.Range("Q19").Value = "$/bandeja 8 kilos"
$ws.Range("S19").Value = 2438
$ws.Range("T19").Value = 8

# Row 20 <- original row 18
$ws.Range("D20").Value = 44505
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 300
$ws.Range("N20").Value = 19000
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 19500
$ws.Range("Q20").Value = "$/bandeja 8 kilos"
$ws.Range("S20").Value = 2438
$ws.Range("T20").Value = 8

